{"js": "// The document opens with the Title, Author and Abstract paragraphs each\n// split across many small runs (one run per word / space, e.g. \"Questions:\",\n// \" \", \"Introduction\", \" \", \"to\", ...). The commit just normalizes those\n// three paragraphs down to a single run each -- the visible wording is\n// unchanged, only the run structure is consolidated. We reproduce that by\n// replacing each paragraph's text with itself (Word/Office.js collapses a\n// replaced range into one run).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// style name -> full (unchanged) paragraph text that should end up in a\n// single run.\nconst targetTextByStyle = {\n  \"Title\": \"Questions: Introduction to sigma notation\",\n  \"Author\": \"Ifan Howells-Baines, Mark Toner\",\n  \"Abstract\": \"Questions relating to the guide on introduction to sigma notation.\",\n};\n\nfor (const paragraph of paragraphs.items) {\n  const targetText = targetTextByStyle[paragraph.style];\n  if (targetText !== undefined && paragraph.text === targetText) {\n    paragraph.insertText(targetText, Word.InsertLocation.replace);\n    delete targetTextByStyle[paragraph.style]; // only touch the first match\n  }\n}\n\nawait context.sync();\n", "ps1": "# The Title, Author and Abstract paragraphs each arrive split across many\n# small runs (one run per word/space, e.g. \"Questions:\", \" \", \"Introduction\",\n# \" \", \"to\", ...). The commit just normalizes those three paragraphs down to\n# a single run each -- the wording itself is unchanged. Running Find/Replace\n# with the exact same text over a paragraph's own range makes Word collapse\n# it to one run, so we do that per style, matched by the current (unchanged)\n# text to make sure we touch the right paragraph.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    @{ Style = \"Title\";    Text = \"Questions: Introduction to sigma notation\" },\n    @{ Style = \"Author\";   Text = \"Ifan Howells-Baines, Mark Toner\" },\n    @{ Style = \"Abstract\"; Text = \"Questions relating to the guide on introduction to sigma notation.\" }\n)\n\nforeach ($target in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        $rng = $p.Range\n        if ($rng.Style.NameLocal -eq $target.Style -and $rng.Text.TrimEnd(\"`r\", \"`n\") -eq $target.Text) {\n            $rng.Find.Execute($target.Text, $false, $false, $false, $false, $false, $true, 1, $false, $target.Text, 2)\n            break\n        }\n    }\n}\n"}
